# This script applies updated odds values to Sheet1, matching the
# source diff (Jogos_da_Semana_FlashScore updated odds/markets).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("J2").Value = 1.14
$ws.Range("K2").Value = 5.5

# Row 3
$ws.Range("K3").Value = 8

# Row 4
$ws.Range("K4").Value = 7.5

# Row 5
$ws.Range("G5").Value = 3.3
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 2
$ws.Range("R5").Value = 1.5
$ws.Range("S5").Value = 2.5
$ws.Range("V5").Value = 12
$ws.Range("X5").Value = 23
$ws.Range("Y5").Value = 26
$ws.Range("AA5").Value = 7.5
$ws.Range("AC5").Value = 34
$ws.Range("AD5").Value = 101
$ws.Range("AE5").Value = 11
$ws.Range("AF5").Value = 12
$ws.Range("AG5").Value = 9
$ws.Range("AH5").Value = 19
$ws.Range("AI5").Value = 15
$ws.Range("AJ5").Value = 21

# Row 6
$ws.Range("N6").Value = 1.67
$ws.Range("O6").Value = 2.15

# Row 7
$ws.Range("G7").Value = 1.95
$ws.Range("H7").Value = 3.75
$ws.Range("I7").Value = 3.6
$ws.Range("K7").Value = 13
$ws.Range("N7").Value = 1.75
$ws.Range("O7").Value = 2.05
$ws.Range("T7").Value = 8.5
$ws.Range("U7").Value = 10
$ws.Range("V7").Value = 8.5
$ws.Range("W7").Value = 17
$ws.Range("X7").Value = 15
$ws.Range("Z7").Value = 13
$ws.Range("AE7").Value = 12
$ws.Range("AF7").Value = 19
$ws.Range("AG7").Value = 12
$ws.Range("AH7").Value = 41
$ws.Range("AI7").Value = 26

# Row 8
$ws.Range("I8").Value = 4

# Row 10
$ws.Range("G10").Value = 1.9
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 4.75
$ws.Range("U10").Value = 7.5
$ws.Range("V10").Value = 9

# Row 11
$ws.Range("N11").Value = 2.15
$ws.Range("O11").Value = 1.67

# Row 13
$ws.Range("G13").Value = 3.7
$ws.Range("H13").Value = 3.25
$ws.Range("I13").Value = 1.91
$ws.Range("U13").Value = 19
$ws.Range("AB13").Value = 19
$ws.Range("AD13").Value = 900
$ws.Range("AE13").Value = 6
$ws.Range("AF13").Value = 8.5
$ws.Range("AG13").Value = 9
$ws.Range("AH13").Value = 17

# Row 14
$ws.Range("J14").Value = 1.05
$ws.Range("K14").Value = 11
$ws.Range("N14").Value = 1.98
$ws.Range("O14").Value = 1.83

# Row 15
$ws.Range("J15").Value = 1.05
$ws.Range("K15").Value = 11

# Row 16
$ws.Range("G16").Value = 19.5
$ws.Range("H16").Value = 6.7
$ws.Range("I16").Value = 1.1
$ws.Range("R16").Value = 2.8
$ws.Range("S16").Value = 1.39
$ws.Range("T16").Value = 45
$ws.Range("U16").Value = 175
$ws.Range("V16").Value = 60
$ws.Range("X16").Value = 350
$ws.Range("Z16").Value = 15
$ws.Range("AA16").Value = 14
$ws.Range("AB16").Value = 37
$ws.Range("AE16").Value = 6.8
$ws.Range("AF16").Value = 4.9
$ws.Range("AH16").Value = 4.9

# Row 17
$ws.Range("G17").Value = 1.88
$ws.Range("H17").Value = 3.65
$ws.Range("I17").Value = 3.3
$ws.Range("T17").Value = 7.4
$ws.Range("U17").Value = 8.5
$ws.Range("X17").Value = 11.75
$ws.Range("Z17").Value = 13
$ws.Range("AA17").Value = 6.4
$ws.Range("AB17").Value = 11.25
$ws.Range("AE17").Value = 10
$ws.Range("AF17").Value = 15.5
$ws.Range("AG17").Value = 9.75
$ws.Range("AH17").Value = 35

# Row 18
$ws.Range("G18").Value = 2.15
$ws.Range("I18").Value = 3.3
$ws.Range("U18").Value = 11
$ws.Range("AE18").Value = 12
